$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column N: "Num de Resolucion" header + first resolution number ---
$ws.Range("N1").Value = "Num de Resolucion"

# Store the resolution number as TEXT (matches how the sibling "Id solicitante"
# column stores numeric-looking identifiers as text), not as a numeric value.
$ws.Range("N2").Value = "'65"
$ws.Range("N2").Style = "Normal"

# --- Row 3's request status moved from PROCESADA back to PENDIENTE ---
$ws.Range("M3").Value = "PENDIENTE"

# --- Column widths for the newly-used L/M columns ---
$ws.Columns("L").ColumnWidth = 19
$ws.Columns("M").ColumnWidth = 15.8

# --- Restore the active selection to M3 like in the edited workbook ---
$ws.Range("M3").Select()
